# Generate Report for Archive
#
# 1. Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all share the same
#    string).
# 2. The status columns narrow to match the shorter text: Overview columns
#    E & F, and the "Status" column (C) on the zh-cn / de-de detail sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status values ---------------------------------------------
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Narrow the status columns ---------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.55
$wsOverview.Columns.Item(6).ColumnWidth = 12.55

$wsZhCn.Columns.Item(3).ColumnWidth = 12.55
$wsDeDe.Columns.Item(3).ColumnWidth = 12.55
